$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Force the cell to be stored as text so numeric-looking strings
    # (e.g. "1.00", "0.000229", "171.90") keep their exact literal form
    # instead of being reinterpreted as numbers by Excel.
    $range.Value = "'" + $text
}


$ws.Range("D2").Value = '43.751.95'
$ws.Range("E2").Value = '  +4.75%  '

$ws.Range("D3").Value = '2.282.16'
$ws.Range("E3").Value = '  +2.43%  '

$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue $ws.Range("D5") '231.65'
$ws.Range("E5").Value = '  +0.12%  '

Set-TextValue $ws.Range("D6") '0.629'
$ws.Range("E6").Value = '  +0.95%  '

Set-TextValue $ws.Range("D7") '64.48'
$ws.Range("E7").Value = '  +6.43%  '

$ws.Range("E8").Value = '  +0.04%  '

Set-TextValue $ws.Range("D9") '0.425'
$ws.Range("E9").Value = '  +5.00%  '

Set-TextValue $ws.Range("D10") '0.0965'
$ws.Range("E10").Value = '  +7.11%  '

Set-TextValue $ws.Range("D11") '57.79'
$ws.Range("E11").Value = '  -0.83%  '

Set-TextValue $ws.Range("D12") '26.39'
$ws.Range("E12").Value = '  +16.11%  '

$ws.Range("E13").Value = '  +0.57%  '

$ws.Range("D14").Value = '2.624.56'
$ws.Range("E14").Value = '  +2.44%  '

Set-TextValue $ws.Range("D15") '15.82'
$ws.Range("E15").Value = '  +1.31%  '

Set-TextValue $ws.Range("D16") '5.93'
$ws.Range("E16").Value = '  +5.50%  '

$ws.Range("E17").Value = '  +2.32%  '

$ws.Range("D18").Value = '2.308.97'

$ws.Range("D19").Value = '43.740.32'
$ws.Range("E19").Value = '  +4.87%  '

$ws.Range("E20").Value = '  +4.51%  '

Set-TextValue $ws.Range("D21") '73.46'
$ws.Range("E21").Value = '  +1.43%  '

$ws.Range("E22").Value = '  +1.40%  '

Set-TextValue $ws.Range("D23") '250.57'
$ws.Range("E23").Value = '  +1.05%  '

$ws.Range("E24").Value = '  +0.15%  '

$ws.Range("E25").Value = '  +7.75%  '

Set-TextValue $ws.Range("D26") '2.33'
$ws.Range("E26").Value = '  +0.99%  '

Set-TextValue $ws.Range("D27") '10.01'
$ws.Range("E27").Value = '  +3.62%  '

Set-TextValue $ws.Range("D28") '171.90'
$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("E29").Value = '  -2.42%  '

Set-TextValue $ws.Range("D30") '20.61'
$ws.Range("E30").Value = '  +3.49%  '

$ws.Range("E31").Value = '  +3.61%  '

$ws.Range("E32").Value = '  +4.77%  '

Set-TextValue $ws.Range("D33") '0.121'
$ws.Range("E33").Value = '  -0.23%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D34") '5.26'
$ws.Range("E34").Value = '  +4.92%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D35") '0.0696'
$ws.Range("E35").Value = '  +6.64%  '

$ws.Range("E36").Value = '  +0.83%  '

Set-TextValue $ws.Range("D37") '6.77'
$ws.Range("E37").Value = '  +3.25%  '

Set-TextValue $ws.Range("D38") '3.79'
$ws.Range("E38").Value = '  +4.82%  '

Set-TextValue $ws.Range("D39") '2.36'
$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("E40").Value = '  +3.70%  '

$ws.Range("B41").Value = 'BinanceUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  +0.12%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D42") '11.16'
$ws.Range("E42").Value = '  +30.08%  '

$ws.Range("B43").Value = 'TerraClassic'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range("D43") '0.000229'
$ws.Range("E43").Value = '  -2.79%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D44") '4.74'
$ws.Range("E44").Value = '  +6.27%  '

$ws.Range("E45").Value = '  -0.78%  '

$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("E47").Value = '  +0.91%  '

Set-TextValue $ws.Range("D48") '98.30'
$ws.Range("E48").Value = '  -0.20%  '

$ws.Range("D49").Value = '1.489.09'
$ws.Range("E49").Value = '  +1.22%  '

Set-TextValue $ws.Range("D50") '16.96'
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("E51").Value = '  +2.52%  '
